$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "n" to the ReasonToReject column (J) for rows whose Approved/Rejected (I) is "Rejected"
$ws.Range("J2").Value = "n"
$ws.Range("J8").Value = "n"
$ws.Range("J10").Value = "n"
$ws.Range("J16").Value = "n"

# Update the active selection to match the saved cursor position
$ws.Range("J16").Select()
